# Updated chemical id mapping
# Populate "Friendly Name" (col B) and "Data type" (col D) for the schema
# terms listed in column A of the environmentalSample sheet, and adjust
# the related column widths / selection to match the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("environmentalSample")

# -- Friendly Name (col B) / Data type (col D) per schema term (rows 2-30) --

$ws.Range("B2").Value  = "Sample Number"

$ws.Range("B3").Value  = "Date sampled"
$ws.Range("D3").Value  = "date"

$ws.Range("D4").Value  = "string"

$ws.Range("B5").Value  = "Technology"
$ws.Range("D5").Value  = "string"

$ws.Range("B6").Value  = "Sample ID"
$ws.Range("D6").Value  = "integer"

$ws.Range("B7").Value  = "ZF Lims Id"
$ws.Range("D7").Value  = "integer"

$ws.Range("B8").Value  = "CAS Number"
$ws.Range("D8").Value  = "string"

$ws.Range("B9").Value  = "Client"
$ws.Range("D9").Value  = "string"

$ws.Range("B10").Value = "Sample Name"
$ws.Range("D10").Value = "string"

$ws.Range("B11").Value = "Latitude"
$ws.Range("D11").Value = "float"

$ws.Range("B12").Value = "Longitude"
$ws.Range("D12").Value = "float"

$ws.Range("B13").Value = "Location Name"
$ws.Range("D13").Value = "string"

$ws.Range("B14").Value = "Location Description"
$ws.Range("D14").Value = "string"

$ws.Range("B15").Value = "Alternate Name"
$ws.Range("D15").Value = "string"

$ws.Range("B16").Value = "Chemical ID"
$ws.Range("D16").Value = "integer"

$ws.Range("B17").Value = "Sample Concentration"
$ws.Range("D17").Value = "float"

$ws.Range("D18").Value = "string"

$ws.Range("B19").Value = "Concentration Unity"
$ws.Range("D19").Value = "string"

$ws.Range("B20").Value = "Sample Molar Concentration"
$ws.Range("D20").Value = "float"

$ws.Range("B21").Value = "Sample Molar Concentration Unit"
$ws.Range("D21").Value = "string"

$ws.Range("B22").Value = "Environmental Concentration"
$ws.Range("D22").Value = "float"

$ws.Range("B24").Value = "Environmental Concentration Unit"

$ws.Range("B25").Value = "Environmental Molar Concentration"

$ws.Range("B26").Value = "Environmental Molar Concentration Unit"

$ws.Range("B27").Value = "Parent Sample Number"

$ws.Range("B28").Value = "Child Sample Number"

$ws.Range("B29").Value = "Project Name"

$ws.Range("B30").Value = "Project Information"

# -- Column width tweaks for the newly-populated Friendly Name / Description columns --
$ws.Columns.Item(2).ColumnWidth = 21.6
$ws.Columns.Item(3).ColumnWidth = 14

# -- Restore the cursor/selection left by the author at B18 --
$ws.Range("B18").Select() | Out-Null
